$d = $word.ActiveDocument

$replacements = @(
    @("879÷6=146, 3", "274÷2=137, 0"),
    @("726÷9=80, 6", "687÷2=343, 1"),
    @("770÷5=154, 0", "777÷7=111, 0"),
    @("417÷8=52, 1", "786÷9=87, 3"),
    @("414÷5=82, 4", "569÷7=81, 2"),
    @("640÷5=128, 0", "876÷3=292, 0"),
    @("580÷2=290, 0", "435÷9=48, 3"),
    @("105÷3=35, 0", "970÷7=138, 4"),
    @("520÷4=130, 0", "742÷4=185, 2"),
    @("834÷9=92, 6", "860÷6=143, 2"),
    @("956÷4=239, 0", "359÷3=119, 2"),
    @("558÷3=186, 0", "979÷6=163, 1"),
    @("542÷3=180, 2", "950÷5=190, 0"),
    @("836÷3=278, 2", "620÷4=155, 0"),
    @("763÷4=190, 3", "585÷4=146, 1"),
    @("522÷9=58, 0", "553÷3=184, 1"),
    @("763÷9=84, 7", "349÷4=87, 1"),
    @("166÷9=18, 4", "568÷2=284, 0"),
    @("696÷3=232, 0", "346÷8=43, 2"),
    @("614÷2=307, 0", "298÷8=37, 2"),
    @("787÷8=98, 3", "808÷8=101, 0"),
    @("434÷6=72, 2", "833÷9=92, 5"),
    @("801÷3=267, 0", "943÷4=235, 3"),
    @("740÷6=123, 2", "961÷4=240, 1"),
    @("412÷6=68, 4", "773÷5=154, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
